$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2 values (B2/C2 swap to Artn/Gfra3 stay same text but shared string index changes internally;
# values for numeric columns change per new TPM data)
$ws.Range("B2").Value = "Artn"
$ws.Range("C2").Value = "Gfra3"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 1.121263666666666
$ws.Range("H2").Value = 3.363791
$ws.Range("I2").Value = 0.8978163344397481
$ws.Range("J2").Value = 0.9294754023256565
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 0.5936135
$ws.Range("N2").Value = 1.187227
$ws.Range("O2").Value = 1
$ws.Range("P2").Value = 1
$ws.Range("Q2").Value = 0.6655972495928332
$ws.Range("R2").Value = 3.993583497557
$ws.Range("S2").Value = 0.8978163344397481
$ws.Range("T2").Value = 0.9294754023256565

# Update row 3 values
$ws.Range("B3").Value = "Artn"
$ws.Range("C3").Value = "Gfra3"
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 0.127615
$ws.Range("H3").Value = 0.25523
$ws.Range("I3").Value = 0.1021836655602519
$ws.Range("J3").Value = 0.07052459767434344
$ws.Range("K3").Value = 2
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 0.5936135
$ws.Range("N3").Value = 1.187227
$ws.Range("O3").Value = 1
$ws.Range("P3").Value = 1
$ws.Range("Q3").Value = 0.07575398680250001
$ws.Range("R3").Value = 0.30301594721
$ws.Range("S3").Value = 0.1021836655602519
$ws.Range("T3").Value = 0.07052459767434344

# Delete row 4 entirely (was Resolving-Mac / Artn / Gfra3 / MuSCs row)
$ws.Rows("4:4").Delete()
